$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 61606.0
$ws.Range("B2").Value = 959.189894067713
$ws.Range("C2").Value = 189555821588.66934
$ws.Range("D2").Value = 42592607.80249771
$ws.Range("E2").Value = 7022.249894095746
# Row 3
$ws.Range("A3").Value = 128310.0
$ws.Range("B3").Value = 959.1539798685804
$ws.Range("C3").Value = 367287043737.24603
$ws.Range("D3").Value = 82623666.2343884
$ws.Range("E3").Value = 13076.153979869561
# Row 4
$ws.Range("A4").Value = 194658.0
$ws.Range("B4").Value = 959.4361003105458
$ws.Range("C4").Value = 544049221222.74066
$ws.Range("D4").Value = 122392377.41744444
$ws.Range("E4").Value = 19134.45610031789
# Row 5
$ws.Range("A5").Value = 261218.0
$ws.Range("B5").Value = 960.386714627813
$ws.Range("C5").Value = 720600517043.9127
$ws.Range("D5").Value = 162186595.43715218
$ws.Range("E5").Value = 25204.226714730376
# Row 6
$ws.Range("A6").Value = 328080.0
$ws.Range("B6").Value = 959.847856083703
$ws.Range("C6").Value = 897302011138.6176
$ws.Range("D6").Value = 201952842.19197965
$ws.Range("E6").Value = 31267.34785630211
# Row 7
$ws.Range("A7").Value = 394220.0
$ws.Range("B7").Value = 960.0436898664883
$ws.Range("C7").Value = 1073935494192.568
$ws.Range("D7").Value = 241697017.07991722
$ws.Range("E7").Value = 37306.24368997436
# Row 8
$ws.Range("A8").Value = 461551.0
$ws.Range("B8").Value = 959.709630921987
$ws.Range("C8").Value = 1252100435729.3523
$ws.Range("D8").Value = 281731345.67952704
$ws.Range("E8").Value = 43388.229630830305
# Row 9
$ws.Range("A9").Value = 527246.0
$ws.Range("B9").Value = 959.3459052917375
$ws.Range("C9").Value = 1427175598612.1926
$ws.Range("D9").Value = 321118580.16557485
$ws.Range("E9").Value = 49407.425904989745
# Row 10
$ws.Range("A10").Value = 592720.0
$ws.Range("B10").Value = 960.396852341065
$ws.Range("C10").Value = 1600496783929.2544
$ws.Range("D10").Value = 360294859.32489014
$ws.Range("E10").Value = 55522.05685181454
# Row 11
$ws.Range("A11").Value = 660392.0
$ws.Range("B11").Value = 958.9296945018265
$ws.Range("C11").Value = 1780673563110.9258
$ws.Range("D11").Value = 400690710.35086316
$ws.Range("E11").Value = 61545.72969375472
# Row 12
$ws.Range("A12").Value = 727264.0
$ws.Range("B12").Value = 959.5025394509189
$ws.Range("C12").Value = 1959575455748.444
$ws.Range("D12").Value = 440887911.98668927
$ws.Range("E12").Value = 67595.80253845913
# Row 13
$ws.Range("A13").Value = 794004.0
$ws.Range("B13").Value = 959.587515006598
$ws.Range("C13").Value = 2136067240961.0576
$ws.Range("D13").Value = 480597338.23787904
$ws.Range("E13").Value = 73663.98751373764
# Row 14
$ws.Range("A14").Value = 860101.0
$ws.Range("B14").Value = 959.2632294303097
$ws.Range("C14").Value = 2311436776124.4585
$ws.Range("D14").Value = 519938093.626502
$ws.Range("E14").Value = 79698.46322788863
# Row 15
$ws.Range("A15").Value = 924358.0
$ws.Range("B15").Value = 960.0883765571574
$ws.Range("C15").Value = 2487459982375.748
$ws.Range("D15").Value = 559379767.6260053
$ws.Range("E15").Value = 85762.40837475332
# Row 16
$ws.Range("A16").Value = 992993.0
$ws.Range("B16").Value = 960.174869023051
$ws.Range("C16").Value = 2663212411581.236
$ws.Range("D16").Value = 599523136.1580825
$ws.Range("E16").Value = 91884.17486695321
# Row 17
$ws.Range("A17").Value = 1059417.0
$ws.Range("B17").Value = 959.3743938405069
$ws.Range("C17").Value = 2842277330265.0845
$ws.Range("D17").Value = 639502286.9609377
$ws.Range("E17").Value = 97866.81439152543
# Row 18
$ws.Range("A18").Value = 1124602.0
$ws.Range("B18").Value = 959.9263814501488
$ws.Range("C18").Value = 3014400957027.0337
$ws.Range("D18").Value = 678268381.422689
$ws.Range("E18").Value = 103960.42637887696
# Row 19
$ws.Range("A19").Value = 1190813.0
$ws.Range("B19").Value = 959.0748038180634
$ws.Range("C19").Value = 3189496492973.313
$ws.Range("D19").Value = 717879419.0625107
$ws.Range("E19").Value = 109968.39480101124
# Row 20
$ws.Range("A20").Value = 1262278.0
$ws.Range("B20").Value = 960.1549265765523
$ws.Range("C20").Value = 3376054708811.765
$ws.Range("D20").Value = 759850985.6109288
$ws.Range("E20").Value = 116061.69492351632
# Row 21
$ws.Range("A21").Value = 1323753.0
$ws.Range("B21").Value = 959.572156967292
$ws.Range("C21").Value = 3542649357210.5303
$ws.Range("D21").Value = 797123077.0861025
$ws.Range("E21").Value = 122082.77215366767
# Row 22
$ws.Range("A22").Value = 1392713.0
$ws.Range("B22").Value = 959.1007444528294
$ws.Range("C22").Value = 3726301845633.7637
$ws.Range("D22").Value = 838480264.9300482
$ws.Range("E22").Value = 128121.78074091373
# Row 23
$ws.Range("A23").Value = 1461131.0
$ws.Range("B23").Value = 959.0397887334505
$ws.Range("C23").Value = 3903805982523.4297
$ws.Range("D23").Value = 878415235.5710466
$ws.Range("E23").Value = 134195.27978497578
# Row 24
$ws.Range("A24").Value = 1523284.0
$ws.Range("B24").Value = 959.8017490610532
$ws.Range("C24").Value = 4072713798388.471
$ws.Range("D24").Value = 916647565.8308747
$ws.Range("E24").Value = 140356.84174511474
# Row 25
$ws.Range("A25").Value = 1594419.0
$ws.Range("B25").Value = 959.6870048339687
$ws.Range("C25").Value = 4257270510865.1167
$ws.Range("D25").Value = 958122208.2267183
$ws.Range("E25").Value = 146293.36700070018
# Row 26
$ws.Range("A26").Value = 1660525.0
$ws.Range("B26").Value = 959.5364761268561
$ws.Range("C26").Value = 4430576889449.647
$ws.Range("D26").Value = 997794979.777978
$ws.Range("E26").Value = 152395.03647179194
# Row 27
$ws.Range("A27").Value = 1725007.0
$ws.Range("B27").Value = 959.5680441268114
$ws.Range("C27").Value = 4604153352767.957
$ws.Range("D27").Value = 1036271944.6848135
$ws.Range("E27").Value = 158450.6880395968
# Row 28
$ws.Range("A28").Value = 1792100.0
$ws.Range("B28").Value = 959.5338035421674
$ws.Range("C28").Value = 4784259317007.264
$ws.Range("D28").Value = 1076760974.179747
$ws.Range("E28").Value = 164556.37379880302
# Row 29
$ws.Range("A29").Value = 1856169.0
$ws.Range("B29").Value = 960.185570107248
$ws.Range("C29").Value = 4958875668344.966
$ws.Range("D29").Value = 1115923259.450139
$ws.Range("E29").Value = 170599.145565162
# Row 30
$ws.Range("A30").Value = 1922952.0
$ws.Range("B30").Value = 959.7761884227557
$ws.Range("C30").Value = 5137701534486.436
$ws.Range("D30").Value = 1155957864.20482
$ws.Range("E30").Value = 176622.25618327453
# Row 31
$ws.Range("A31").Value = 1988909.0
$ws.Range("B31").Value = 959.2856940001129
$ws.Range("C31").Value = 5317504613427.499
$ws.Range("D31").Value = 1196170990.564489
$ws.Range("E31").Value = 182678.88568864582
